$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 140837
$ws.Range("E2").Value = 3302
$ws.Range("F2").Value = 3302
$ws.Range("G2").Value = 2907
$ws.Range("H2").Value = 2349
$ws.Range("I2").Value = 2350
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 286104
$ws.Range("L2").Value = 262381
$ws.Range("M2").Value = 23723
$ws.Range("N2").Value = 23723
$ws.Range("P2").Value = 447
$ws.Range("Q2").Value = 19343
$ws.Range("R2").Value = -18533
$ws.Range("S2").Value = 340
$ws.Range("T2").Value = 184
$ws.Range("V2").Value = 792
$ws.Range("W2").Value = 2.34
$ws.Range("X2").Value = 1.67
$ws.Range("Y2").Value = 10.66
$ws.Range("Z2").Value = 0.88
$ws.Range("AA2").Value = 1106.01
$ws.Range("AB2").Value = 5252.05
$ws.Range("AC2").Value = 2628
$ws.Range("AD2").Value = 9.890000000000001
$ws.Range("AE2").Value = 29501
$ws.Range("AF2").Value = 0.88
$ws.Range("AG2").Value = 750
$ws.Range("AH2").Value = 2.88
$ws.Range("AI2").Value = 25.67
$ws.Range("AJ2").Value = 89400000
$ws.Range("O2").ClearContents()
$ws.Range("U2").ClearContents()

# Row 3
$ws.Range("D3").Value = 153444
$ws.Range("E3").Value = 3310
$ws.Range("F3").Value = 3310
$ws.Range("G3").Value = 2976
$ws.Range("H3").Value = 2123
$ws.Range("I3").Value = 2123
$ws.Range("K3").Value = 327814
$ws.Range("L3").Value = 302787
$ws.Range("M3").Value = 25027
$ws.Range("N3").Value = 25027
$ws.Range("P3").Value = 447
$ws.Range("Q3").Value = 14978
$ws.Range("R3").Value = -17114
$ws.Range("S3").Value = 3608
$ws.Range("T3").Value = 481
$ws.Range("V3").Value = 5161
$ws.Range("W3").Value = 2.16
$ws.Range("X3").Value = 1.38
$ws.Range("Y3").Value = 8.710000000000001
$ws.Range("Z3").Value = 0.6899999999999999
$ws.Range("AA3").Value = 1209.84
$ws.Range("AB3").Value = 5582.76
$ws.Range("AC3").Value = 2375
$ws.Range("AD3").Value = 15.18
$ws.Range("AE3").Value = 31387
$ws.Range("AF3").Value = 1.15
$ws.Range("AG3").Value = 750
$ws.Range("AH3").Value = 2.08
$ws.Range("AI3").Value = 28.17
$ws.Range("AJ3").Value = 89400000
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("U3").ClearContents()

# Row 4
$ws.Range("D4").Value = 153484
$ws.Range("E4").Value = 5439
$ws.Range("F4").Value = 5439
$ws.Range("G4").Value = 5266
$ws.Range("H4").Value = 4099
$ws.Range("I4").Value = 4099
$ws.Range("K4").Value = 370574
$ws.Range("L4").Value = 342092
$ws.Range("M4").Value = 28482
$ws.Range("N4").Value = 28482
$ws.Range("P4").Value = 447
$ws.Range("Q4").Value = 21399
$ws.Range("R4").Value = -24051
$ws.Range("S4").Value = 1494
$ws.Range("T4").Value = 833
$ws.Range("V4").Value = 6040
$ws.Range("W4").Value = 3.54
$ws.Range("X4").Value = 2.67
$ws.Range("Y4").Value = 15.32
$ws.Range("Z4").Value = 1.17
$ws.Range("AA4").Value = 1201.1
$ws.Range("AB4").Value = 6355.63
$ws.Range("AC4").Value = 4585
$ws.Range("AD4").Value = 6.87
$ws.Range("AE4").Value = 35720
$ws.Range("AF4").Value = 0.88
$ws.Range("AG4").Value = 1350
$ws.Range("AH4").Value = 4.29
$ws.Range("AI4").Value = 26.26
$ws.Range("AJ4").Value = 89400000
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("U4").ClearContents()

# Row 5
$ws.Range("D5").Value = 158868
$ws.Range("E5").Value = 6308
$ws.Range("F5").Value = 6308
$ws.Range("G5").Value = 6219
$ws.Range("H5").Value = 4644
$ws.Range("I5").Value = 4644
$ws.Range("K5").Value = 404925
$ws.Range("L5").Value = 373888
$ws.Range("M5").Value = 31037
$ws.Range("N5").Value = 31037
$ws.Range("P5").Value = 447
$ws.Range("Q5").Value = 17334
$ws.Range("R5").Value = -23776
$ws.Range("S5").Value = 3877
$ws.Range("T5").Value = 322
$ws.Range("V5").Value = 11168
$ws.Range("W5").Value = 3.97
$ws.Range("X5").Value = 2.92
$ws.Range("Y5").Value = 15.61
$ws.Range("Z5").Value = 1.2
$ws.Range("AA5").Value = 1204.64
$ws.Range("AB5").Value = 6927.36
$ws.Range("AC5").Value = 5195
$ws.Range("AD5").Value = 9.050000000000001
$ws.Range("AE5").Value = 38925
$ws.Range("AF5").Value = 1.21
$ws.Range("AG5").Value = 1500
$ws.Range("AH5").Value = 3.19
$ws.Range("AI5").Value = 25.75
$ws.Range("AJ5").Value = 89400000
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("U5").ClearContents()

# Row 6
$ws.Range("D6").Value = 157466
$ws.Range("E6").Value = 5335
$ws.Range("F6").Value = 5335
$ws.Range("G6").Value = 5111
$ws.Range("H6").Value = 3735
$ws.Range("I6").Value = 3735
$ws.Range("K6").Value = 440922
$ws.Range("L6").Value = 399608
$ws.Range("M6").Value = 41314
$ws.Range("N6").Value = 41314
$ws.Range("P6").Value = 447
$ws.Range("Q6").Value = 18983
$ws.Range("R6").Value = -21614
$ws.Range("S6").Value = 3757
$ws.Range("T6").Value = 474
$ws.Range("V6").Value = 11171
$ws.Range("W6").Value = 3.39
$ws.Range("X6").Value = 2.37
$ws.Range("Y6").Value = 10.32
$ws.Range("Z6").Value = 0.88
$ws.Range("AA6").Value = 967.26
$ws.Range("AB6").Value = 9226.290000000001
$ws.Range("AC6").Value = 4178
$ws.Range("AD6").Value = 9.83
$ws.Range("AE6").Value = 51813
$ws.Range("AF6").Value = 0.79
$ws.Range("AG6").Value = 1130
$ws.Range("AH6").Value = 2.75
$ws.Range("AI6").Value = 24.12
$ws.Range("AJ6").Value = 89400000
$ws.Range("U6").ClearContents()

# Row 7
$ws.Range("D7").Value = 121846
$ws.Range("E7").Value = 4105
$ws.Range("G7").Value = 3820
$ws.Range("H7").Value = 2692
$ws.Range("I7").Value = 2715
$ws.Range("K7").Value = 460170
$ws.Range("L7").Value = 414433
$ws.Range("M7").Value = 45737
$ws.Range("N7").Value = 45640
$ws.Range("P7").Value = 447
$ws.Range("W7").Value = 3.37
$ws.Range("X7").Value = 2.21
$ws.Range("Y7").Value = 6.25
$ws.Range("Z7").Value = 0.6
$ws.Range("AA7").Value = 906.12
$ws.Range("AC7").Value = 3037
$ws.Range("AD7").Value = 7.24
$ws.Range("AE7").Value = 57239
$ws.Range("AF7").Value = 0.38
$ws.Range("AG7").Value = 906
$ws.Range("AH7").Value = 4.12
$ws.Range("AI7").Value = 29.84
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()

# Row 8
$ws.Range("D8").Value = 126954
$ws.Range("E8").Value = 4084
$ws.Range("G8").Value = 3940
$ws.Range("H8").Value = 2825
$ws.Range("I8").Value = 2855
$ws.Range("K8").Value = 479554
$ws.Range("L8").Value = 431775
$ws.Range("M8").Value = 47779
$ws.Range("N8").Value = 47740
$ws.Range("P8").Value = 447
$ws.Range("W8").Value = 3.22
$ws.Range("X8").Value = 2.23
$ws.Range("Y8").Value = 6.12
$ws.Range("Z8").Value = 0.6
$ws.Range("AA8").Value = 903.6900000000001
$ws.Range("AC8").Value = 3194
$ws.Range("AD8").Value = 6.89
$ws.Range("AE8").Value = 59873
$ws.Range("AF8").Value = 0.37
$ws.Range("AG8").Value = 950
$ws.Range("AH8").Value = 4.32
$ws.Range("AI8").Value = 29.74
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()

# Row 9
$ws.Range("D9").Value = 127168
$ws.Range("E9").Value = 4951
$ws.Range("G9").Value = 4738
$ws.Range("H9").Value = 3405
$ws.Range("I9").Value = 3430
$ws.Range("K9").Value = 497186
$ws.Range("L9").Value = 446946
$ws.Range("M9").Value = 50241
$ws.Range("N9").Value = 50350
$ws.Range("P9").Value = 447
$ws.Range("W9").Value = 3.89
$ws.Range("X9").Value = 2.68
$ws.Range("Y9").Value = 6.99
$ws.Range("Z9").Value = 0.7
$ws.Range("AA9").Value = 889.61
$ws.Range("AC9").Value = 3837
$ws.Range("AD9").Value = 5.73
$ws.Range("AE9").Value = 63146
$ws.Range("AF9").Value = 0.35
$ws.Range("AG9").Value = 1094
$ws.Range("AH9").Value = 4.97
$ws.Range("AI9").Value = 28.52
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
